# Natmi following Dr Hou advice
# Recompute the Btc -> Erbb3 ligand-receptor table for every combination of
# Sending cluster x Target cluster across {FAPs, sCs}, expanding the result
# from 3 rows (FAPs -> M1/M2/sCs) to the full 8-row cross product
# ({FAPs,sCs} x {FAPs,M1,M2,sCs}) with refreshed statistics.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
#          D Target cluster, E..T the various expression/specificity metrics.
$data = @(
    @("FAPs", "Btc", "Erbb3", "FAPs", 3, 1, 0.637502, 1.912506, 0.6183852387299037, 0.6183852387299037, 2, 0.6666666666666666, 0.2048153333333333, 0.614446, 0.06238051989214988, 0.06238051989214988, 0.1305701846306667, 1.175131661676, 0.03857519268560261, 0.03857519268560261),
    @("FAPs", "Btc", "Erbb3", "M1", 3, 1, 0.637502, 1.912506, 0.6183852387299037, 0.6183852387299037, 3, 1, 0.7041406666666666, 2.112422, 0.2144598265618379, 0.2144598265618379, 0.4488910832813333, 4.040019749532, 0.1326187910464159, 0.1326187910464159),
    @("FAPs", "Btc", "Erbb3", "M2", 3, 1, 0.637502, 1.912506, 0.6183852387299037, 0.6183852387299037, 3, 1, 0.7300083333333333, 2.190025, 0.2223383309140356, 0.2223383309140356, 0.4653817725166666, 4.18843595265, 0.1374907418410843, 0.1374907418410843),
    @("FAPs", "Btc", "Erbb3", "sCs", 3, 1, 0.637502, 1.912506, 0.6183852387299037, 0.6183852387299037, 3, 1, 1.644357666666667, 4.933073, 0.5008213226319767, 0.5008213226319767, 1.048281301215333, 9.434531710938002, 0.309700513156801, 0.309700513156801),
    @("sCs", "Btc", "Erbb3", "FAPs", 2, 0.6666666666666666, 0.393412, 1.180236, 0.3816147612700962, 0.3816147612700962, 2, 0.6666666666666666, 0.2048153333333333, 0.614446, 0.06238051989214988, 0.06238051989214988, 0.08057680991733335, 0.7251912892560001, 0.02380532720654727, 0.02380532720654727),
    @("sCs", "Btc", "Erbb3", "M1", 2, 0.6666666666666666, 0.393412, 1.180236, 0.3816147612700962, 0.3816147612700962, 3, 1, 0.7041406666666666, 2.112422, 0.2144598265618379, 0.2144598265618379, 0.2770173879546667, 2.493156491592, 0.081841035515422, 0.081841035515422),
    @("sCs", "Btc", "Erbb3", "M2", 2, 0.6666666666666666, 0.393412, 1.180236, 0.3816147612700962, 0.3816147612700962, 3, 1, 0.7300083333333333, 2.190025, 0.2223383309140356, 0.2223383309140356, 0.2871940384333334, 2.5847463459, 0.08484758907295137, 0.08484758907295137),
    @("sCs", "Btc", "Erbb3", "sCs", 2, 0.6666666666666666, 0.393412, 1.180236, 0.3816147612700962, 0.3816147612700962, 3, 1, 1.644357666666667, 4.933073, 0.5008213226319767, 0.5008213226319767, 0.6469100383586668, 5.822190345228001, 0.1911208094751756, 0.1911208094751756)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (the header) is untouched; data rows start at row 2 and now run
# through row 9 (dimension grows from A1:T4 to A1:T9).
$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
